$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46060 -> 46061) for every data row (rows 2 through 23).
$ws.Range("C2:C23").Value = 46061
